# storeCharts.xlsx edit:
#  - Update the three existing chart titles to consolidated Finnish text.
#  - Insert a new "Sheet4" worksheet (between "Sheet3" and "Sheet2") holding
#    the LRS / ZRS redundancy timing data.
#  - Add a clustered bar chart on the new sheet comparing LRS vs ZRS across
#    the four service tiers, with a value-axis title and bottom legend,
#    mirroring the look of the other charts in the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update existing chart titles
# ---------------------------------------------------------------------------

# Sheet1 -> "Azure Files" chart
$ws1 = $wb.Worksheets.Item("Sheet1")
$chart1 = $ws1.ChartObjects().Item(1).Chart
$chart1.HasTitle = $true
$chart1.ChartTitle.Text = "Azure Files -palvelun eri palvelumallien kirjoitusnopeudet tiedostoille"

# Sheet3 -> "Blob Storage" chart
$ws3 = $wb.Worksheets.Item("Sheet3")
$chart2 = $ws3.ChartObjects().Item(1).Chart
$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "Blob Storage -palvelun eri palvelumallien kirjoitusnopeudet tiedostoille"

# Sheet2 -> "Managed Disk" chart
$ws2 = $wb.Worksheets.Item("Sheet2")
$chart3 = $ws2.ChartObjects().Item(1).Chart
$chart3.HasTitle = $true
$chart3.ChartTitle.Text = "Managed Disk -palvelun eri palvelumallien kirjoitusnopeudet tiedostoille"

# ---------------------------------------------------------------------------
# 2. Insert the new "Sheet4" worksheet (lands before the ActiveSheet, i.e.
#    between "Sheet3" and "Sheet2", matching the workbook's original
#    activeTab / ActiveSheet = "Sheet2").
# ---------------------------------------------------------------------------

$ws4 = $wb.Worksheets.Add()
$ws4.Name = "Sheet4"

$ws4.Range("A1").Value = "Korkealuokkainen"
$ws4.Range("A2").Value = "Tapahtumaoptimoitu"
$ws4.Range("A3").Value = "Kuuma"
$ws4.Range("A4").Value = "Viileä"

$ws4.Range("B1").Value = 2685
$ws4.Range("B2").Value = 2721
$ws4.Range("B3").Value = 2716
$ws4.Range("B4").Value = 3559

$ws4.Range("C1").Value = 2689
$ws4.Range("C2").Value = 2687
$ws4.Range("C3").Value = 2686
$ws4.Range("C4").Value = 2685

# ---------------------------------------------------------------------------
# 3. Build the LRS/ZRS bar chart on the new sheet
# ---------------------------------------------------------------------------

$co4 = $ws4.ChartObjects().Add(620, 170, 545, 255)
$chart4 = $co4.Chart
$chart4.ChartType = 57   # xlBarClustered
$chart4.SetSourceData($ws4.Range("A1:C4"))

$cg4 = $chart4.ChartGroups(1)
$cg4.GapWidth = 182
$chart4.VaryColors = $false

$s1 = $chart4.SeriesCollection().Item(1)
$s1.Name = "LRS"
$s2 = $chart4.SeriesCollection().Item(2)
$s2.Name = "ZRS"

$chart4.HasTitle = $true
$chart4.ChartTitle.Text = "Azure Files -palvelun 300 MB:n tiedoston tallentamisnopeudet eri redundanssi-alueilla"

$valAx4 = $chart4.Axes(2, 1)
$valAx4.HasTitle = $true
$valAx4.AxisTitle.Text = "aika (ms)"

$chart4.HasLegend = $true
$chart4.Legend.Position = -4107  # xlLegendPositionBottom

# ---------------------------------------------------------------------------
# 4. Minor selection/cursor nudges to mirror the author's click-through
# ---------------------------------------------------------------------------

$ws1.Activate()
$ws1.Range("A1:B4").Select()

$ws3.Activate()
$ws3.Range("F19").Select()

$ws4.Activate()
$ws4.Range("F25").Select()

$ws2.Activate()
$ws2.Range("M35").Select()

# Leave "Sheet4" as the active sheet, matching the saved workbook state.
$ws4.Activate()
